$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 130 (pushes old rows 130-137 down to 132-139),
# copying formatting (incl. the date number format on column D) from the row above.
$ws.Rows.Item(130).Insert()
$ws.Rows.Item(130).Insert()

# New row 130: week of 2021-11-16, "1a (cosecha lavada)" quality, Provincia de Melipilla
$ws.Cells.Item(130, 1).Value2 = 11
$ws.Cells.Item(130, 2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item(130, 3).Value2 = "Bíobío"
$ws.Cells.Item(130, 4).Value2 = 44516
$ws.Cells.Item(130, 5).Value2 = 8
$ws.Cells.Item(130, 6).Value2 = 100114001
$ws.Cells.Item(130, 7).Value2 = "Papa"
$ws.Cells.Item(130, 8).Value2 = "Asterix"
$ws.Cells.Item(130, 9).Value2 = "1a (cosecha lavada)"
$ws.Cells.Item(130, 10).Value2 = 450
$ws.Cells.Item(130, 11).Value2 = 11000
$ws.Cells.Item(130, 12).Value2 = 12000
$ws.Cells.Item(130, 13).Value2 = 11444
$ws.Cells.Item(130, 14).Value2 = "`$/malla 25 kilos"
$ws.Cells.Item(130, 15).Value2 = "Provincia de Melipilla"
$ws.Cells.Item(130, 16).Value2 = 458
$ws.Cells.Item(130, 17).Value2 = 25
$ws.Cells.Item(130, 18).Value2 = "Hortaliza"

# New row 131: week of 2021-11-16, "1a nueva(o)" quality, Provincia de Melipilla
$ws.Cells.Item(131, 1).Value2 = 11
$ws.Cells.Item(131, 2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item(131, 3).Value2 = "Bíobío"
$ws.Cells.Item(131, 4).Value2 = 44516
$ws.Cells.Item(131, 5).Value2 = 8
$ws.Cells.Item(131, 6).Value2 = 100114001
$ws.Cells.Item(131, 7).Value2 = "Papa"
$ws.Cells.Item(131, 8).Value2 = "Asterix"
$ws.Cells.Item(131, 9).Value2 = "1a nueva(o)"
$ws.Cells.Item(131, 10).Value2 = 350
$ws.Cells.Item(131, 11).Value2 = 10000
$ws.Cells.Item(131, 12).Value2 = 11000
$ws.Cells.Item(131, 13).Value2 = 10429
$ws.Cells.Item(131, 14).Value2 = "`$/saco 25 kilos"
$ws.Cells.Item(131, 15).Value2 = "Provincia de Melipilla"
$ws.Cells.Item(131, 16).Value2 = 417
$ws.Cells.Item(131, 17).Value2 = 25
$ws.Cells.Item(131, 18).Value2 = "Hortaliza"
